# "added 4wk low sales check"
# Update the forecasted sales figures (and their dependent Inventory Coverage /
# Seasonality Index columns) on the "Forecast Comparison" sheet, then refresh
# the derived aggregate metrics on the "Summary" sheet to match.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row -> [MyForecast (D), Inventory Coverage (H), Seasonality Index (L)]
$rows = @{
    2  = @(184, 16.29, 1.01)
    3  = @(183, 15.37, 1.16)
    4  = @(169, 15.56, 0.92)
    5  = @(165, 14.92, 1.11)
    6  = @(171, 13.43, 0.99)
    7  = @(180, 11.81, 0.93)
    8  = @(177, 10.99, 0.95)
    9  = @(168, 10.52, 0.99)
    10 = @(166, 9.64,  1.02)
    11 = @(169, 8.49,  1.02)
    12 = @(181, 6.99,  1)
    13 = @(176, 6.16,  0.8)
    14 = @(154, 5.9,   1.16)
    15 = @(138, 5.46,  1.14)
    16 = @(144, 4.28,  1.12)
    17 = @(150, 3.15,  1.01)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $wsForecast.Range("D$r").Value = $vals[0]
    $wsForecast.Range("H$r").Value = $vals[1]
    $wsForecast.Range("L$r").Value = $vals[2]
}

# Refresh the Summary sheet's derived totals / extremes from the new
# MyForecast (column D) figures. These cells hold their numbers as TEXT
# (same as the original file). Prefix the entry with an apostrophe so Excel
# stores it as text instead of auto-converting to a number, then strip the
# resulting "quote prefix" formatting so the cell style is left untouched.
$summaryCells = @("B9", "B10", "B11", "B12", "B14")
$summaryValues = @("2675", "1397", "701", "184", "138")
for ($i = 0; $i -lt $summaryCells.Length; $i++) {
    $cell = $wsSummary.Range($summaryCells[$i])
    $cell.Value = "'" + $summaryValues[$i]
    $cell.ClearFormats()
}
